$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New grid of marker names for columns A and B (rows 1-91).
# Column C, and any rows beyond 91, are no longer used.
$data = @(
    @(1, "Rajagopal 2015", "Santos 2017"),
    @(2, "RACR", "R.Shoulder"),
    @(3, "LACR", "L.Shoulder"),
    @(4, "C7", "C7"),
    @(5, "T2", "T2"),
    @(6, "T7", "T7"),
    @(7, "L1", "L1"),
    @(8, "L3", "L3"),
    @(9, "L5", "L5"),
    @(10, "IJ", "IJ"),
    @(11, "CLAV", ""),
    @(12, "RASH", ""),
    @(13, "RPSH", ""),
    @(14, "LASH", ""),
    @(15, "LPSH", ""),
    @(16, "RSJC", ""),
    @(17, "RUA1", ""),
    @(18, "RUA2", ""),
    @(19, "RUA3", ""),
    @(20, "RLEL", ""),
    @(21, "RMEL", ""),
    @(22, "RFAsuperior", ""),
    @(23, "RFAradius", ""),
    @(24, "RFAulna", ""),
    @(25, "LSJC", ""),
    @(26, "LUA1", ""),
    @(27, "LUA2", ""),
    @(28, "LUA3", ""),
    @(29, "LLEL", ""),
    @(30, "LMEL", ""),
    @(31, "LFAsuperior", ""),
    @(32, "LFAradius", ""),
    @(33, "LFAulna", ""),
    @(34, "RASI", "R.ASIS"),
    @(35, "LASI", "L.ASIS"),
    @(36, "RPSI", "R.PSIS"),
    @(37, "LPSI", "L.PSIS"),
    @(38, "LHJC", ""),
    @(39, "RHJC", ""),
    @(40, "RTH1", ""),
    @(41, "RTH2", ""),
    @(42, "RTH3", ""),
    @(43, "RLFC", "R.Knee"),
    @(44, "RMFC", "R.Knee.Medialv"),
    @(45, "RKJC", ""),
    @(46, "RTB1", ""),
    @(47, "RTB2", ""),
    @(48, "RTB3", ""),
    @(49, "RLMAL", "R.Ankle"),
    @(50, "RMMAL", "R.Ankle.Medialv"),
    @(51, "RAJC", ""),
    @(52, "RCAL", "R.Heel"),
    @(53, "RTOE", "R.MT1"),
    @(54, "RMT5", "R.MT5"),
    @(55, "RMT2", "R.MT2"),
    @(56, "LTH1", ""),
    @(57, "LTH2", ""),
    @(58, "LTH3", ""),
    @(59, "LLFC", "L.Knee"),
    @(60, "LMFC", "L.Knee.Medialv"),
    @(61, "LKCJ", ""),
    @(62, "LTB1", ""),
    @(63, "LTB2", ""),
    @(64, "LTB3", ""),
    @(65, "LLMAL", "L.Ankle"),
    @(66, "LMMAL", "L.Ankle.Medialv"),
    @(67, "LAJC", ""),
    @(68, "LCAL", "L.Heel"),
    @(69, "LTOE", "L.MT1"),
    @(70, "LMT5", "L.MT5"),
    @(71, "LMT2", "L.MT2"),
    @(72, "REJC", ""),
    @(73, "LEJC", ""),
    @(74, "R_tibial_plateau", ""),
    @(75, "L_tibial_plateau", ""),
    @(76, "T10", ""),
    @(77, "NAVE", ""),
    @(78, "XYPH", "PX"),
    @(79, "RGTRO", "R.GTR"),
    @(80, "LGTRO", "L.GTR"),
    @(81, "RFAX", "R.HF"),
    @(82, "LFAX", "L.HF"),
    @(83, "RTTC", "R.TT"),
    @(84, "LTTC", "L.TT"),
    @(85, "", "COG"),
    @(86, "", "R.Front.Head"),
    @(87, "", "L.Front.Head"),
    @(88, "", "R.Back.Head"),
    @(89, "", "L.Back.Head"),
    @(90, "", "L.Iliac.Crestv"),
    @(91, "", "R.Iliac.Crestv")
)

# Clear the whole previously used range (A1:C83) first so stale column C
# data and any leftover cells beyond the new extent are removed.
$ws.Range("A1:C83").Clear()

# Rows/cells that hold brand-new marker names (not previously present
# anywhere in the workbook). These must be entered in a specific order
# (matching how the source file was authored) so the resulting shared
# string table lines up: RMT2, LMT2, RTTC, LTTC, RFAX, LFAX.
$deferredCells = @{
    55 = $true
    71 = $true
    81 = $true
    82 = $true
    83 = $true
    84 = $true
}

foreach ($row in $data) {
    $r = $row[0]
    $aVal = $row[1]
    $bVal = $row[2]

    if ($aVal -ne "" -and -not $deferredCells.ContainsKey($r)) {
        $ws.Cells.Item($r, 1).Value = $aVal
    }
    if ($bVal -ne "") {
        $ws.Cells.Item($r, 2).Value = $bVal
    }
}

# Now enter the brand-new marker names in the order they were first typed
# by the original author, so new shared-string entries are appended in
# that same order: RMT2, LMT2, RTTC, LTTC, RFAX, LFAX.
$ws.Cells.Item(55, 1).Value = "RMT2"
$ws.Cells.Item(71, 1).Value = "LMT2"
$ws.Cells.Item(83, 1).Value = "RTTC"
$ws.Cells.Item(84, 1).Value = "LTTC"
$ws.Cells.Item(81, 1).Value = "RFAX"
$ws.Cells.Item(82, 1).Value = "LFAX"

$ws.Range("C5").Select()
